# Update the "专家发言" (expert-panel) prompt cell F12 with the new, expanded
# prompt text, resize the row to fit the longer content, and move the
# selection/cursor to F12 (matching the author's edit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @"
专家发言/模拟一场会议，提供10个专家的解答.###思考方式:**增强AI的复杂推理和思维链能力，解决高难度问题。
问题分解：
明确问题的关键要素，逐步将其分解为更小、更易处理的子问题。
目标是识别问题中的不同层次，并为每个子问题找到突破口。
思维链 (Chain of Thought - CoT)：
对每个子问题构建一系列清晰的思维步骤，记录并解释每个推理过程。
确保思维链条完整、严谨，避免遗漏任何关键逻辑环节。
假设验证：

对每个可能的解决方案进行假设测试，分析其合理性和潜在局限性。
根据反馈修正错误，确保假设的准确性和可行性。
多策略思考：

如果现有方法未能得出理想结果，迅速切换策略，尝试不同路径，探索备选方案。
综合多个角度分析问题，形成更全面的解决方案。
推理总结：

汇总所有推理步骤，结合多种思考路径，最终得出一个具有说服力的答案。
阐明每个结论背后的逻辑及推理依据，确保解释条理清晰、严密。###

Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3

"@

$ws.Range("F12").Value = $newText

$ws.Rows.Item(12).RowHeight = 409.6

[void]$ws.Range("F12").Select()
